# Refresh cryptos list: updated prices / 1h volume deltas, and re-ranked a
# few coins whose order changed (rows 45-47: dogwifhat / Fetch.AI / ApeXProtocol).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.016.39"
$ws.Range("E2").Value = "  +4.11%  "

$ws.Range("D3").Value = "3.460.81"
$ws.Range("E3").Value = "  +3.77%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'584.69"
$ws.Range("E5").Value = "  +5.74%  "

$ws.Range("D6").Value = "'187.13"
$ws.Range("E6").Value = "  +8.15%  "

$ws.Range("D7").Value = "'0.633"
$ws.Range("E7").Value = "  +1.16%  "

$ws.Range("D8").Value = "3.454.21"
$ws.Range("E8").Value = "  +3.86%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "'0.172"
$ws.Range("E10").Value = "  +0.38%  "

$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("D12").Value = "'56.24"
$ws.Range("E12").Value = "  +5.47%  "

$ws.Range("E13").Value = "  +0.39%  "

$ws.Range("E14").Value = "  +3.51%  "

$ws.Range("D15").Value = "4.020.12"
$ws.Range("E15").Value = "  +4.38%  "

$ws.Range("D16").Value = "'18.71"
$ws.Range("E16").Value = "  +3.28%  "

$ws.Range("D17").Value = "3.461.61"
$ws.Range("E17").Value = "  +4.04%  "

$ws.Range("D18").Value = "66.959.68"
$ws.Range("E18").Value = "  +4.22%  "

$ws.Range("D19").Value = "'12.11"
$ws.Range("E19").Value = "  +3.39%  "

$ws.Range("E20").Value = "  -2.15%  "

$ws.Range("D21").Value = "'1.02"
$ws.Range("E21").Value = "  +3.20%  "

$ws.Range("D22").Value = "'490.84"
$ws.Range("E22").Value = "  +9.56%  "

$ws.Range("E23").Value = "  +7.87%  "

$ws.Range("D24").Value = "'16.56"
$ws.Range("E24").Value = "  +20.97%  "

$ws.Range("D25").Value = "'4.44"
$ws.Range("E25").Value = "  +9.61%  "

$ws.Range("D26").Value = "'89.68"
$ws.Range("E26").Value = "  +3.51%  "

$ws.Range("E27").Value = "  +2.36%  "

$ws.Range("D28").Value = "'10.89"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("D29").Value = "'9.04"
$ws.Range("E29").Value = "  +5.28%  "

$ws.Range("D30").Value = "'31.31"
$ws.Range("E30").Value = "  +1.52%  "

$ws.Range("E31").Value = "  +9.37%  "

$ws.Range("D32").Value = "'599.59"
$ws.Range("E32").Value = "  +5.38%  "

$ws.Range("D33").Value = "'11.69"
$ws.Range("E33").Value = "  +2.73%  "

$ws.Range("D34").Value = "'63.82"
$ws.Range("E34").Value = "  +2.60%  "

$ws.Range("E35").Value = "  +4.49%  "

$ws.Range("E36").Value = "  +6.57%  "

$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("D38").Value = "'36.46"
$ws.Range("E38").Value = "  +3.40%  "

$ws.Range("D39").Value = "'3.54"

$ws.Range("E40").Value = "  +4.50%  "

$ws.Range("E41").Value = "  +3.25%  "

$ws.Range("D42").Value = "3.234.38"
$ws.Range("E42").Value = "  +5.71%  "

$ws.Range("E43").Value = "  +6.29%  "

$ws.Range("E44").Value = "  +3.32%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "'2.78"
$ws.Range("E45").Value = "  +23.74%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.50"
$ws.Range("E46").Value = "  +2.83%  "

$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "'3.23"
$ws.Range("E47").Value = "  +2.44%  "

$ws.Range("E48").Value = "  +1.03%  "

$ws.Range("D49").Value = "'3.28"
$ws.Range("E49").Value = "  +13.23%  "

$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").Value = "'8.71"
$ws.Range("E51").Value = "  +6.46%  "
